# Duplicate the last slide (the cover/title slide with picture) and move
# the new copy to the very beginning of the deck.
$p = $ppt.ActivePresentation

$sourceSlide = $p.Slides.Item($p.Slides.Count)
$newSlideRange = $sourceSlide.Duplicate()
$newSlide = $newSlideRange.Item(1)
$newSlide.MoveTo(1)
